# The Jostens QA automation workbook records the outcome of each
# automated test-suite run by writing a "<Status> - <timestamp>" string
# into a fixed "result" cell on each worksheet (Login, School Search,
# Product Search, Shopping Cart, Checkout, Payment). This script updates
# those result cells to reflect the latest test-suite run.

$wb = $excel.ActiveWorkbook

$results = @(
    @{ Sheet = "Login";           Cell = "G2"; Value = "Success - 2020/12/19 12:01:53" }
    @{ Sheet = "Login";           Cell = "G3"; Value = "Success - 2020/12/19 12:01:56" }
    @{ Sheet = "School Search";   Cell = "C2"; Value = "Success - 2020/12/19 12:01:59" }
    @{ Sheet = "School Search";   Cell = "C3"; Value = "Success - 2020/12/19 12:02:01" }
    @{ Sheet = "Product Search";  Cell = "K2"; Value = "Success - 2020/12/19 12:02:27" }
    @{ Sheet = "Product Search";  Cell = "K3"; Value = "Success - 2020/12/19 12:02:48" }
    @{ Sheet = "Product Search";  Cell = "K4"; Value = "Success - 2020/12/19 12:03:09" }
    @{ Sheet = "Shopping Cart";   Cell = "G2"; Value = "Success - 2020/12/19 12:03:11" }
    @{ Sheet = "Shopping Cart";   Cell = "G3"; Value = "Success - 2020/12/19 12:03:11" }
    @{ Sheet = "Shopping Cart";   Cell = "G4"; Value = "Success - 2020/12/19 12:03:11" }
    @{ Sheet = "Checkout";        Cell = "P2"; Value = "Success - 2020/12/19 12:03:21" }
    @{ Sheet = "Checkout";        Cell = "P3"; Value = "Success - 2020/12/19 12:03:32" }
    @{ Sheet = "Checkout";        Cell = "P4"; Value = "Success - 2020/12/19 12:03:40" }
    @{ Sheet = "Payment";         Cell = "C2"; Value = "Success - 2020/12/19 12:03:50" }
)

foreach ($result in $results) {
    $ws = $wb.Worksheets.Item($result.Sheet)
    $ws.Range($result.Cell).Value = $result.Value
}

# The "Payment" sheet's result cell (C2) carries the same quote-prefixed
# cell style as its row-mates (A2, D2). Plain Value assignment resets a
# cell to the default style, so re-apply the original formatting by
# copying it from a sibling cell that already has it.
$paymentSheet = $wb.Worksheets.Item("Payment")
$styleSource = $paymentSheet.Range("A2")
$styleSource.Copy()
$paymentSheet.Range("C2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false
